# Automatische test-sync: 2025-08-28 21:08:50
#
# Adds the new "Opvolging bestelling" log entry (2025-08-28 21:08:16) to the
# Logs sheet, bumps the Dashboard "Overig" tally row, and extends the chart
# series / conditional-formatting ranges so they keep covering the sheets'
# used ranges.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet - append row 22 with the new mail log entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A22").Value = "Opvolging bestelling"
$logs.Range("B22").Value = "mailmind.test@zohomail.eu"
$logs.Range("D22").Value = "Overig"
$logs.Range("F22").Value = "2025-08-28 21:08:16"
$logs.Range("G22").Value = "Nee"
$logs.Range("H22").Value = "Ja"
$logs.Range("I22").Value = "Nee"
$logs.Range("J22").Value = "Nee"

# Extend the existing conditional-formatting blocks (D/G/H/I/J) from row 21
# down to row 22 while keeping their dxf/priority/formula assignments intact.
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range("$($col)2:$($col)21")
    $newRange = $logs.Range("$($col)2:$($col)22")
    $rules = $oldRange.FormatConditions
    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet - add the "Overig" tally row
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Overig"
$dash.Range("B4").Value = 1

# ---------------------------------------------------------------------
# 3. Chart - widen the category/value series ranges to include row 4
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$4,Dashboard!`$B`$2:`$B`$4,1)"
